# Remove the (redundant/default) <w:contextualSpacing w:val="0"/> element
# from every paragraph's <w:pPr> in the document body.
#
# There is no exposed ParagraphFormat.ContextualSpacing property on this
# COM surface, so we round-trip the body through its raw WordOpenXML:
# read it, strip the element via a targeted regex, and write it back with
# Range.InsertXML (the only supported way to replace a range's underlying
# OOXML content - Range.WordOpenXML / Document.WordOpenXML are read-only).

$d = $word.ActiveDocument
$body = $d.Content

$xml = $body.WordOpenXML

# Drop every self-closing <w:contextualSpacing .../> tag, regardless of
# attribute order/spacing, anywhere in the package part.
$pattern = '<w:contextualSpacing\b[^/>]*/>'
$newXml = [System.Text.RegularExpressions.Regex]::Replace($xml, $pattern, '')

if ($newXml -ne $xml) {
    $body.InsertXML($newXml) | Out-Null
}
